$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142, shifting the existing rows 142-151 down to 143-152.
$ws.Rows(142).Insert()

# Populate the newly inserted row 142 with the new weekly data point.
$ws.Cells.Item(142, 1).Value  = 9
$ws.Cells.Item(142, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(142, 3).Value  = "Metropolitana"
$ws.Cells.Item(142, 4).Value  = 44461
$ws.Cells.Item(142, 5).Value  = 13
$ws.Cells.Item(142, 6).Value  = 300000001
$ws.Cells.Item(142, 7).Value  = "Rabanito"
$ws.Cells.Item(142, 8).Value  = "Sin especificar"
$ws.Cells.Item(142, 9).Value  = "Primera"
$ws.Cells.Item(142, 10).Value = 7900
$ws.Cells.Item(142, 11).Value = 3500
$ws.Cells.Item(142, 12).Value = 4000
$ws.Cells.Item(142, 13).Value = 3747
$ws.Cells.Item(142, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(142, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(142, 16).Value = 37
$ws.Cells.Item(142, 17).Value = 100
$ws.Cells.Item(142, 18).Value = "Hortaliza"
